# "updated main GSC export data"
#
# The GSC export window rolled forward by one day on the "Chart"
# sheet:
#   - the oldest date (2025-11-06, row 2) drops off the front
#   - every remaining row shifts up by one row
#   - a new date (2026-02-03) is appended as the new last row, with
#     the "Valid" count the export reported for it (28 — same as the
#     previous last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest row (row 2, date 2025-11-06) and shift rows 3:90 up
# to become rows 2:89. xlShiftUp = -4162.
$ws.Range("A2:C2").Delete(-4162)

# Append the new row for the new date at the bottom of the (now
# 89-row) data block.
$newRow = 90

$dateCell = $ws.Cells.Item($newRow, 1)
# Force the cell to text first so the date string is stored verbatim
# (as the rest of column A is) instead of being auto-converted into a
# date serial number, then drop the format override so the cell keeps
# the sheet's plain default style.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-02-03"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 28
